$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.093.74"
$ws.Range("E2").Value = "  -3.08%  "
$ws.Range("D3").Value = "3.501.31"
$ws.Range("E3").Value = "  -5.11%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.69"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.88%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.06"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -3.46%  "
$ws.Range("D7").Value = "3.493.96"
$ws.Range("E7").Value = "  -5.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.605"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.35%  "
$ws.Range("E9").Value = "  -0.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.188"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -5.34%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.47"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.55%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.581"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -4.42%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "46.47"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -4.97%  "
$ws.Range("E14").Value = "  -4.54%  "
$ws.Range("D15").Value = "4.058.31"
$ws.Range("E15").Value = "  -5.43%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.51"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -5.05%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "617.33"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -9.09%  "
$ws.Range("D18").Value = "68.999.92"
$ws.Range("E18").Value = "  -3.38%  "
$ws.Range("D19").Value = "3.496.26"
$ws.Range("E19").Value = "  -5.53%  "
$ws.Range("E20").Value = "  -0.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.31"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -3.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.10"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -3.75%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.882"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -6.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "15.85"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -8.80%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "97.05"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -4.97%  "
$ws.Range("E26").Value = "  -4.78%  "
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("E28").Value = "  -6.69%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.29"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -9.62%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.57"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -6.98%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.15"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -7.85%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.48"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -6.95%  "
$ws.Range("E33").Value = "  -8.32%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.96"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -5.00%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "632.91"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +7.41%  "
$ws.Range("E36").Value = "  -4.14%  "
$ws.Range("E37").Value = "  -5.33%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.43"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -15.60%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "56.50"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -4.26%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0447"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.81%  "
$ws.Range("E42").Value = "  -5.83%  "
$ws.Range("D43").Value = "3.356.34"
$ws.Range("E43").Value = "  -8.65%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.326"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -6.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "32.70"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -7.09%  "
$ws.Range("D46").Value = "0.0₃0689"
$ws.Range("E46").Value = "  -9.85%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.56"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -7.40%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.77"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -3.31%  "
$ws.Range("E49").Value = "  -2.56%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "132.33"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.97%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.60"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +13.83%  "
